$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("T2A")
$ws.Delete()
